$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: variant labels gain a space ("Вариант1" -> "Вариант 1", ...)
#    A1/B1/G1 keep their text ("Раздел" / "Вопрос" / "Правильный ответ").
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Вариант 1"
$ws.Range("D1").Value = "Вариант 2"
$ws.Range("E1").Value = "Вариант 3"
$ws.Range("F1").Value = "Вариант 4"

# ---------------------------------------------------------------------------
# 2. Column A used to hold the plain section number (1-5) with style index 1
#    (numFmtId 1, i.e. "0"). It now holds a text label ("N секция") with the
#    default (unstyled) format, so clear the explicit style before writing
#    the new text - copy the style from a cell that is already unstyled.
# ---------------------------------------------------------------------------
$defaultStyle = $ws.Range("B1").Style
$ws.Range("A2:A24").Style = $defaultStyle

$ws.Range("A2").Value = "1 секция"
$ws.Range("A3").Value = "1 секция"
$ws.Range("A4").Value = "1 секция"
$ws.Range("A5").Value = "1 секция"
$ws.Range("A6").Value = "1 секция"

$ws.Range("A7").Value = "2 секция"
$ws.Range("A8").Value = "2 секция"
$ws.Range("A9").Value = "2 секция"
$ws.Range("A10").Value = "2 секция"
$ws.Range("A11").Value = "2 секция"

$ws.Range("A12").Value = "3 секция"
$ws.Range("A13").Value = "3 секция"
$ws.Range("A14").Value = "3 секция"
$ws.Range("A15").Value = "3 секция"

$ws.Range("A16").Value = "4 секция"
$ws.Range("A17").Value = "4 секция"
$ws.Range("A18").Value = "4 секция"
$ws.Range("A19").Value = "4 секция"
$ws.Range("A20").Value = "4 секция"

$ws.Range("A21").Value = "5 секция"
$ws.Range("A22").Value = "5 секция"
$ws.Range("A23").Value = "5 секция"
$ws.Range("A24").Value = "5 секция"

# ---------------------------------------------------------------------------
# 3. Column G ("Правильный ответ") previously always duplicated column F
#    (the last option). Fix it to point at the actually-correct option for
#    the rows whose correct answer is NOT the last option.
# ---------------------------------------------------------------------------
$ws.Range("G6").Value = 206
$ws.Range("G8").Value = "Бейсбол"
$ws.Range("G9").Value = "Возвращение Лесси"
$ws.Range("G10").Value = "Роуэн Аткинсон"
$ws.Range("G11").Value = 5
$ws.Range("G12").Value = "Джелато"
$ws.Range("G14").Value = "Гренландия"
$ws.Range("G15").Value = "Луара"
$ws.Range("G16").Value = "Зубы"
$ws.Range("G17").Value = "Голова"
$ws.Range("G18").Value = "Котята"
$ws.Range("G20").Value = 100
$ws.Range("G21").Value = "Октоторп"
$ws.Range("G22").Value = "J"
$ws.Range("G24").Value = "Аллодоксафобия"

# ---------------------------------------------------------------------------
# 4. Column A is narrower now (27.28515625 -> 15.85546875 character-width
#    units). ColumnWidth is quantized to 1/6-character steps by this host,
#    so feed it the character width net of the 5/6 padding constant and let
#    it round to the nearest reachable value.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 15.0221354166667

# ---------------------------------------------------------------------------
# 5. View tweaks: drop the custom zoom (back to 100%) and move the
#    selection to G24.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("G24").Select()
